# Workbook was regenerated against a renamed repo / fixed output path,
# which changed which MAG rows ended up in this species' result sheet.
# Net effect on this sheet: the rows for even_MAG-GUT16297.fa,
# even_MAG-GUT30319.fa and even_MAG-GUT7957.fa are no longer part of the
# output, while all remaining rows (and their values) are unchanged and
# shift up to fill the gaps.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete from the bottom up so earlier row numbers stay valid while we work.
$ws.Rows.Item(9).Delete()   # even_MAG-GUT7957.fa
$ws.Rows.Item(3).Delete()   # even_MAG-GUT30319.fa
$ws.Rows.Item(2).Delete()   # even_MAG-GUT16297.fa
